# feat: add 2022-Q1 data
#
# The workbook tracks quarterly fund-holding snapshots, one sheet per quarter,
# plus a trailing "总计" (grand total) summary sheet. Adding a new quarter
# means:
#   1. A new "2022-Q1" sheet (same shape as the other quarterly sheets) is
#      inserted right before "总计".
#   2. The "总计" sheet gets a new first data row summarizing 2022-Q1, with
#      the existing rows pushed down and their running index renumbered.
#
# To get the sheetId/relationship-id bookkeeping to come out the same way a
# real Excel "insert sheet" + "rename" sequence would, we repurpose the
# existing "总计" sheet object (which keeps its original identity) as the new
# "2022-Q1" sheet, and create a brand-new sheet for the (now one-row-bigger)
# "总计" table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: assign a value that must be stored as text (not auto-converted to a
# number/date by Excel), while stripping the left-over "quote prefix" style
# so the cell keeps the default (no explicit) style afterwards.
# ---------------------------------------------------------------------------
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

$oldTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------------
# Step 1: spin up the new "总计" sheet right after the current one, carrying
# over the current totals table (data + number formats + page setup) before
# the source sheet gets reused for 2022-Q1 data.
# ---------------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $oldTotal)
$oldTotal.Range("A1:D6").Copy($newTotal.Range("A1"))
$newTotal.Range("A1").ClearContents()

$newTotal.PageSetup.LeftMargin = 54
$newTotal.PageSetup.RightMargin = 54
$newTotal.PageSetup.TopMargin = 72
$newTotal.PageSetup.BottomMargin = 72
$newTotal.PageSetup.HeaderMargin = 36
$newTotal.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------------
# Step 2: insert the new 2022-Q1 summary row (becomes row 2), push the
# existing rows down, and renumber the running index column (A).
# ---------------------------------------------------------------------------
$newTotal.Rows.Item(2).Insert()

Set-TextValue $newTotal.Range("B2") "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("C2").Style = "Normal"
$newTotal.Range("D2").Value = 0.05
$newTotal.Range("D2").Style = "Normal"

# Row 3 (old row 2) already carries the correct "index" column style (s=2);
# clone it onto the new row 2, then fix up every row's index value.
$newTotal.Range("A3").Copy($newTotal.Range("A2"))
$newTotal.Range("A2").Value = 0
$newTotal.Range("A3").Value = 1
$newTotal.Range("A4").Value = 2
$newTotal.Range("A5").Value = 3
$newTotal.Range("A6").Value = 4
$newTotal.Range("A7").Value = 5

# ---------------------------------------------------------------------------
# Step 3: repurpose the original sheet object as "2022-Q1" (frees up the
# "总计" name) and fill it with this quarter's fund-holding figures, cloning
# the layout of the previous quarter (2021-Q4).
# ---------------------------------------------------------------------------
$oldTotal.Cells.Clear()
$oldTotal.Name = "2022-Q1"

$prevQuarter = $wb.Worksheets.Item(5)
$prevQuarter.Range("A1:H3").Copy($oldTotal.Range("A1"))
$oldTotal.Range("A1").ClearContents()

# Row 2: 001092 / 广发纳斯达克生物科技指数(QDII)（人民币）
Set-TextValue $oldTotal.Range("B2") "001092"
Set-TextValue $oldTotal.Range("D2") "1.34"
Set-TextValue $oldTotal.Range("E2") "82.00"
Set-TextValue $oldTotal.Range("F2") "2.01"
Set-TextValue $oldTotal.Range("G2") "0.0269"
$oldTotal.Range("H2").Value = 9

# Row 3: 001093 / 广发纳斯达克生物科技指数(QDII)（美元）
Set-TextValue $oldTotal.Range("B3") "001093"
Set-TextValue $oldTotal.Range("D3") "1.34"
Set-TextValue $oldTotal.Range("E3") "82.00"
Set-TextValue $oldTotal.Range("F3") "2.01"
Set-TextValue $oldTotal.Range("G3") "0.0269"
$oldTotal.Range("H3").Value = 9

# ---------------------------------------------------------------------------
# Step 4: now "总计" is free, claim it for the sheet built in step 1.
# ---------------------------------------------------------------------------
$newTotal.Name = "总计"
